# Remove the two data rows (RM 232 and SC 92) from the "missing_data" table.
# Row 26 = "RM 232" and the original row 28 = "SC 92"; once row 26 is removed
# the old row 28 shifts up to row 27, so we delete that position next.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# After the deletions the remaining rows shift up two positions. Two of the
# "missing" (blank) data cells in column D move along with their rows, and
# the diff shows the blank/filled status of column D differs for the rows
# that now land on 29 and 33, so set those explicitly to match the target.
$ws.Range("D29").Value = ""
$ws.Range("D33").Value = -14.1
